$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns I and J (headers "I0" and "IF") mirror the styling of the
# existing header row (bold, centered, bordered) by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4
